{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Requisitos\" Heading2 paragraph and the requirement bullet\n// paragraph that follows it (the last two paragraphs of the document) and\n// remove both, restoring the document to its pre-\"Requisitos section\" state.\nconst toDelete = [];\nfor (const p of paragraphs.items) {\n  const text = (p.text || \"\").trim();\n  if (text === \"Requisitos\" || text.indexOf(\"LOQ4237\") !== -1) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the trailing \"Requisitos\" section: the Heading2 paragraph titled\n# \"Requisitos\" and the ListBullet paragraph right after it that lists the\n# weak prerequisite (LOQ4237). Walk paragraphs back-to-front so deleting one\n# doesn't disturb the index of paragraphs still pending a check.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text\n    if ($text.StartsWith(\"Requisitos\") -or $text.Contains(\"LOQ4237\")) {\n        $para.Range.Delete()\n    }\n}\n"}
